# Update the crypto price/volume table (columns D and E, rows 2-51) with
# freshly scraped values. Column D prices are stored as TEXT in the sheet
# (e.g. "65.734.62", "7.50") so numeric-looking values are written with a
# leading apostrophe to force Excel to keep them as text instead of
# coercing them to numbers (which would drop the formatting/trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.734.62'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '2.700.63'
$ws.Range("E3").Value = '  +1.86%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''605.86'
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("D6").Value = '''157.70'
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("E9").Value = '  +5.41%  '
$ws.Range("D10").Value = '''6.04'
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").Value = '''30.20'
$ws.Range("E13").Value = '  +4.32%  '
$ws.Range("E14").Value = '  +10.27%  '
$ws.Range("D15").Value = '3.184.11'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").Value = '65.616.02'
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("D17").Value = '2.700.73'
$ws.Range("E17").Value = '  +2.59%  '
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").Value = '''4.86'
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").Value = '''360.02'
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").Value = '''7.50'
$ws.Range("E21").Value = '  +2.83%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = '''70.48'
$ws.Range("E23").Value = '  +3.49%  '
$ws.Range("E24").Value = '  +3.82%  '
$ws.Range("E25").Value = '  +12.55%  '
$ws.Range("E26").Value = '  -5.23%  '
$ws.Range("E27").Value = '  +3.25%  '
$ws.Range("E28").Value = '  +4.20%  '
$ws.Range("D29").Value = '''8.34'
$ws.Range("E29").Value = '  +2.26%  '
$ws.Range("E30").Value = '  +3.68%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").Value = '''540.21'
$ws.Range("E32").Value = '  +4.07%  '
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").Value = '''6.74'
$ws.Range("E34").Value = '  +6.37%  '
$ws.Range("E35").Value = '  -4.21%  '
$ws.Range("D36").Value = '''0.431'
$ws.Range("E36").Value = '  +1.20%  '
$ws.Range("D37").Value = '''20.72'
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("D38").Value = '''162.59'
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = '''42.82'
$ws.Range("E42").Value = '  +1.50%  '
$ws.Range("D43").Value = '''167.73'
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("D44").Value = '''4.20'
$ws.Range("E44").Value = '  +2.36%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = '''23.54'
$ws.Range("E46").Value = '  +2.94%  '
$ws.Range("E47").Value = '  +2.84%  '
$ws.Range("D48").Value = '''0.661'
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").Value = '''0.0267'
$ws.Range("E49").Value = '  +4.65%  '
$ws.Range("D50").Value = '''21.20'
$ws.Range("E50").Value = '  +9.08%  '
$ws.Range("D51").Value = '''0.0985'
$ws.Range("E51").Value = '  +0.06%  '
